$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the paragraph that currently ends with "Add icon, title." - this is
# the paragraph that (in the original document) also carries the _GoBack
# bookmark at its very end, and after which the new "Expanded" widget notes
# must be inserted.
# ---------------------------------------------------------------------------
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Add icon, title.`r") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not locate the 'Add icon, title.' paragraph"
}

$rPrXml = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-RunXml($text) {
    return '<w:r>' + $rPrXml + '<w:t>' + $text + '</w:t></w:r>'
}

function New-ParaXml($ilvl, $innerRunsXml) {
    return '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="' + $ilvl + '"/><w:numId w:val="2"/></w:numPr>' + $rPrXml + '</w:pPr>' + $innerRunsXml + '</w:p>'
}

# Build up the run-level XML for each of the six new paragraphs first.
$runsImage = New-RunXml "Image:"

$runsAssetImage = '<w:proofErr w:type="spellStart"/>'
$runsAssetImage += New-RunXml "AssetImage"
$runsAssetImage += '<w:proofErr w:type="spellEnd"/>'
$runsAssetImage += New-RunXml ": Image Provider."

$runsExpanded = New-RunXml "Expanded"
$runsMakeChild = New-RunXml "Make child of Row, Column, or Flex expand to full the available space in the main axis."
$runsIfMultiple = New-RunXml "If multiple children are expanded, the available space is divided among them according to the flex factor."
$runsMustBe = New-RunXml "Must be a descendant of a Row, Colum, or flex"

# Now wrap each in its paragraph (with the proper list indent level).
$para1 = New-ParaXml 2 $runsImage
$para2 = New-ParaXml 3 $runsAssetImage
$para3 = New-ParaXml 2 $runsExpanded
$para4 = New-ParaXml 3 $runsMakeChild
$para5 = New-ParaXml 3 $runsIfMultiple
$para6 = New-ParaXml 3 $runsMustBe

$newParas = @($para1, $para2, $para3, $para4, $para5, $para6)

$r = $d.Paragraphs.Item($anchorIndex).Range
$r.Collapse(0)

foreach ($paraXml in $newParas) {
    $r.InsertParagraphAfter()
    $r.Collapse(0)
    $anchorIndex = $anchorIndex + 1
    $target = $d.Paragraphs.Item($anchorIndex)
    $target.Range.InsertXML($paraXml)
    $r = $d.Paragraphs.Item($anchorIndex).Range
    $r.Collapse(0)
}

# ---------------------------------------------------------------------------
# Move the _GoBack bookmark from the "Add icon, title." paragraph to the
# last (empty) paragraph of the document, right before the section break.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
